$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Octubre de 2020 a las 10:27"

# Apply updated country stats / re-sorted rows (Casos totales refresh)
# Row 7
$ws.Cells.Item(7, 2).Value = 1194643
$ws.Cells.Item(7, 3).Value = 9412
$ws.Cells.Item(7, 4).Value = 970296
$ws.Cells.Item(7, 5).Value = 203270
$ws.Cells.Item(7, 7).Value = 186
$ws.Cells.Item(7, 8).Value = 21077

# Row 23
$ws.Cells.Item(23, 2).Value = 316678
$ws.Cells.Item(23, 3).Value = 2611
$ws.Cells.Item(23, 4).Value = 254617
$ws.Cells.Item(23, 5).Value = 56445
$ws.Cells.Item(23, 7).Value = 56
$ws.Cells.Item(23, 8).Value = 5616

# Row 27
$ws.Cells.Item(27, 2).Value = 256071
$ws.Cells.Item(27, 3).Value = 2581
$ws.Cells.Item(27, 4).Value = 183500
$ws.Cells.Item(27, 5).Value = 70942
$ws.Cells.Item(27, 7).Value = 7
$ws.Cells.Item(27, 8).Value = 1629

# Row 59
$ws.Cells.Item(59, 2).Value = 57794
$ws.Cells.Item(59, 3).Value = 10
$ws.Cells.Item(59, 5).Value = 255

# Row 64
$ws.Cells.Item(64, 2).Value = 51382
$ws.Cells.Item(64, 3).Value = 532
$ws.Cells.Item(64, 4).Value = 44406
$ws.Cells.Item(64, 5).Value = 6010
$ws.Cells.Item(64, 7).Value = 3
$ws.Cells.Item(64, 8).Value = 966

# Row 86
$ws.Cells.Item(86, 1).Value = "Tunez"
$ws.Cells.Item(86, 2).Value = 19721
$ws.Cells.Item(86, 3).Value = 1308
$ws.Cells.Item(86, 4).Value = 5032
$ws.Cells.Item(86, 5).Value = 14418
$ws.Cells.Item(86, 7).Value = 6
$ws.Cells.Item(86, 8).Value = 271

# Row 87
$ws.Cells.Item(87, 1).Value = "Grecia"
$ws.Cells.Item(87, 2).Value = 18886
$ws.Cells.Item(87, 4).Value = 9989
$ws.Cells.Item(87, 5).Value = 8504
$ws.Cells.Item(87, 8).Value = 393

# Row 89
$ws.Cells.Item(89, 2).Value = 17160
$ws.Cells.Item(89, 3).Value = 333
$ws.Cells.Item(89, 4).Value = 15423
$ws.Cells.Item(89, 5).Value = 1446
$ws.Cells.Item(89, 7).Value = 7
$ws.Cells.Item(89, 8).Value = 291

# Row 98
$ws.Cells.Item(98, 1).Value = "Eslovaquia"
$ws.Cells.Item(98, 2).Value = 11617
$ws.Cells.Item(98, 3).Value = 679
$ws.Cells.Item(98, 4).Value = 4756
$ws.Cells.Item(98, 5).Value = 6807
$ws.Cells.Item(98, 7).Value = 6
$ws.Cells.Item(98, 8).Value = 54

# Row 99
$ws.Cells.Item(99, 1).Value = "Malasia"
$ws.Cells.Item(99, 2).Value = 11484
$ws.Cells.Item(99, 4).Value = 10014
$ws.Cells.Item(99, 5).Value = 1334
$ws.Cells.Item(99, 8).Value = 136

# Row 100
$ws.Cells.Item(100, 1).Value = "Namibia"
$ws.Cells.Item(100, 2).Value = 11373
$ws.Cells.Item(100, 4).Value = 9083
$ws.Cells.Item(100, 5).Value = 2167
$ws.Cells.Item(100, 8).Value = 123

# Row 101
$ws.Cells.Item(101, 1).Value = "Montenegro"
$ws.Cells.Item(101, 2).Value = 10987
$ws.Cells.Item(101, 4).Value = 7397
$ws.Cells.Item(101, 5).Value = 3420
$ws.Cells.Item(101, 8).Value = 170

# Row 115
$ws.Cells.Item(115, 5).Value = 3333
$ws.Cells.Item(115, 7).Value = 5
$ws.Cells.Item(115, 8).Value = 46

# Row 129
$ws.Cells.Item(129, 1).Value = "Lituania"
$ws.Cells.Item(129, 2).Value = 4956
$ws.Cells.Item(129, 3).Value = 172
$ws.Cells.Item(129, 4).Value = 2466
$ws.Cells.Item(129, 5).Value = 2397
$ws.Cells.Item(129, 7).Value = 1
$ws.Cells.Item(129, 8).Value = 93

# Row 130
$ws.Cells.Item(130, 1).Value = "Surinam"
$ws.Cells.Item(130, 2).Value = 4891
$ws.Cells.Item(130, 4).Value = 4702
$ws.Cells.Item(130, 5).Value = 84
$ws.Cells.Item(130, 8).Value = 105

# Row 131
$ws.Cells.Item(131, 1).Value = "Ruanda"
$ws.Cells.Item(131, 2).Value = 4843
$ws.Cells.Item(131, 4).Value = 3181
$ws.Cells.Item(131, 5).Value = 1633
$ws.Cells.Item(131, 8).Value = 29

# Row 132
$ws.Cells.Item(132, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(132, 2).Value = 4829
$ws.Cells.Item(132, 4).Value = 1914
$ws.Cells.Item(132, 5).Value = 2853
$ws.Cells.Item(132, 8).Value = 62

# Row 142
$ws.Cells.Item(142, 2).Value = 3507
$ws.Cells.Item(142, 3).Value = 57
$ws.Cells.Item(142, 4).Value = 2675
$ws.Cells.Item(142, 5).Value = 767

# Row 207
$ws.Cells.Item(207, 1).Value = "Santa Lucia"

# Row 208
$ws.Cells.Item(208, 1).Value = "Nueva Caledonia"
